# Generate Report for Handoff
# Refresh the "Latest Handoff"/"Latest HO Xliff Generate Date" timestamps for the
# 0a1136ba-86bd-4c30-94ce-37d9e5e864d0 row (row 5) across the Overview, zh-cn and
# de-de sheets to reflect the newly generated handoff report.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G5").Value = "2016-08-17 00:39:20"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H5").Value = "2016-08-17 00:39:15"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H5").Value = "2016-08-17 00:39:20"
